$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for "Espárragos" at row 50.
# This pushes the existing rows 50:86 down one row to become rows 51:87.
$ws.Rows("50:50").Insert()

# Populate the newly inserted row 50 with the new record's data.
$ws.Range("A50").Value = 10
$ws.Range("B50").Value = 'Vega Modelo de Temuco'
$ws.Range("C50").Value = 'La Araucanía'
$ws.Range("D50").Value = 45159
$ws.Range("D50").NumberFormat = $ws.Range("D51").NumberFormat
$ws.Range("E50").Value = 9
$ws.Range("F50").Value = 300000000
$ws.Range("G50").Value = 'Espárragos'
$ws.Range("H50").Value = 'Sin especificar'
$ws.Range("I50").Value = 'Primera'
$ws.Range("J50").Value = 900
$ws.Range("K50").Value = 3500
$ws.Range("L50").Value = 3500
$ws.Range("M50").Value = 3500
$ws.Range("N50").Value = '$/kilo'
$ws.Range("O50").Value = 'Región del Maule'
$ws.Range("P50").Value = 3500
$ws.Range("Q50").Value = 1
$ws.Range("R50").Value = 'Hortaliza'
